# Update "想去人数" (F column) values for rows 5-34 on sheets "展览" and "全部类型".
# Mapping of row -> new value (only rows that actually changed per the diff).

$updates = @{
    5  = 80
    6  = 248
    7  = 34
    8  = 496
    9  = 38
    10 = 1922
    11 = 56
    12 = 72
    13 = 3995
    14 = 32
    15 = 262
    16 = 93
    17 = 69
    18 = 9
    20 = 2708
    22 = 369
    23 = 14
    25 = 53
    27 = 42
    29 = 5
    30 = 33
    32 = 139
    33 = 1568
    34 = 208
}

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
